$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 50, shifting existing rows 50-99 down to 51-100
$ws.Rows.Item(50).Insert()

# Populate the newly inserted row 50 with the new data record
$ws.Range("A50").Value = 10
$ws.Range("B50").Value = "Vega Modelo de Temuco"
$ws.Range("C50").Value = "La Araucanía"
$ws.Range("D50").Value = 45090
$ws.Range("E50").Value = 9
$ws.Range("F50").Value = "Fruta"
$ws.Range("G50").Value = 100108
$ws.Range("H50").Value = "Tropicales y subtropicales"
$ws.Range("I50").Value = 100108003
$ws.Range("J50").Value = "Maracuyá"
$ws.Range("K50").Value = "Sin especificar"
$ws.Range("L50").Value = "Primera"
$ws.Range("M50").Value = 40
$ws.Range("N50").Value = 50000
$ws.Range("O50").Value = 50000
$ws.Range("P50").Value = 50000
$ws.Range("Q50").Value = "$/caja 18 kilos"
$ws.Range("R50").Value = "Región de Arica y Parinacota"
$ws.Range("S50").Value = 2778
$ws.Range("T50").Value = 18
